{"js": "// Word JavaScript API (Office.js) script.\n// Body of: async (context) => { ... }\n//\n// 1) Paragraph \"FL can be categorized ... use cases for FL.\" was split\n//    across two runs (\"..\" + \".\"); normalize it back into a single run\n//    holding the full sentence (no visible text change, just a run merge).\n// 2) After the Federated Transfer Learning (FTL) paragraph, insert two new\n//    paragraphs: one with the new \"Terms like FL and DML ...\" sentence, and\n//    one empty paragraph (both using the same BodyText style/format as their\n//    neighbours), while leaving the existing trailing empty paragraphs\n//    untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- Step 1: merge the two runs that make up the \"FL can be categorized...\" paragraph ---\nconst targetPrefix = \"FL can be categorized according to the distribution of the data held by the clients participating in the modelling.\";\nlet categorizedPara = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(targetPrefix) === 0) {\n    categorizedPara = p;\n    break;\n  }\n}\n\nif (categorizedPara) {\n  const fullText = categorizedPara.text;\n  categorizedPara.getRange().insertText(fullText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Step 2: insert the two new paragraphs after the FTL paragraph ---\nconst ftlPrefix = \"Federated Transfer Learning (FTL) is applicable\";\nlet ftlPara = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(ftlPrefix) === 0) {\n    ftlPara = p;\n    break;\n  }\n}\n\nif (ftlPara) {\n  const ooxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n    '<w:p><w:pPr><w:pStyle w:val=\"BodyText\"/><w:spacing w:line=\"360\" w:lineRule=\"auto\"/><w:ind w:left=\"540\" w:right=\"332\" w:firstLine=\"470\"/></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Terms like FL and Distributed Machine Learning (DML) can create confusion due to their similarities. The main difference lies in the training process: in FL, there is a central server that aggregates updates sent by the clients, whereas in DML, there is no central server; instead, data is spread across different nodes and computations are shared among these nodes. Table 2.4 will help clarify the intricacies of each concept.</w:t></w:r>' +\n    '</w:p>' +\n    '<w:p><w:pPr><w:pStyle w:val=\"BodyText\"/><w:spacing w:line=\"360\" w:lineRule=\"auto\"/><w:ind w:left=\"540\" w:right=\"332\" w:firstLine=\"470\"/></w:pPr></w:p>' +\n    '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n\n  ftlPara.getRange(\"End\").insertOoxml(ooxml, Word.InsertLocation.after);\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is the open document.\n#\n# 1) Paragraph \"FL can be categorized ... use cases for FL.\" was split\n#    across two runs (\"..\" + \".\"); normalize it back into a single run\n#    holding the full sentence (no visible text change, just a run merge).\n# 2) After the Federated Transfer Learning (FTL) paragraph, insert two new\n#    paragraphs: one with the new \"Terms like FL and DML ...\" sentence, and\n#    one empty paragraph (both using the same BodyText style/format as their\n#    neighbours), while leaving the existing trailing empty paragraphs\n#    untouched.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: merge the two runs that make up the \"FL can be categorized...\" paragraph ---\n$targetPrefix = \"FL can be categorized according to the distribution of the data held by the clients participating in the modelling.\"\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.StartsWith($targetPrefix)) {\n        $full = $p.Range\n        $fullText = $full.Text\n        $fullText = $fullText.Substring(0, $fullText.Length - 1)  # drop trailing paragraph mark\n        $r = $d.Range($full.Start, $full.End - 1)\n        $r.Delete()\n        $r2 = $d.Range($full.Start, $full.Start)\n        $r2.InsertBefore($fullText)\n        break\n    }\n}\n\n# --- Step 2: insert the two new paragraphs after the FTL paragraph ---\n$ftlPrefix = \"Federated Transfer Learning (FTL) is applicable\"\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.StartsWith($ftlPrefix)) {\n        $endPos = $p.Range.End\n        $r = $d.Range($endPos, $endPos)\n        $xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:pStyle w:val=\"BodyText\"/><w:spacing w:line=\"360\" w:lineRule=\"auto\"/><w:ind w:left=\"540\" w:right=\"332\" w:firstLine=\"470\"/></w:pPr><w:r><w:t>Terms like FL and Distributed Machine Learning (DML) can create confusion due to their similarities. The main difference lies in the training process: in FL, there is a central server that aggregates updates sent by the clients, whereas in DML, there is no central server; instead, data is spread across different nodes and computations are shared among these nodes. Table 2.4 will help clarify the intricacies of each concept.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"BodyText\"/><w:spacing w:line=\"360\" w:lineRule=\"auto\"/><w:ind w:left=\"540\" w:right=\"332\" w:firstLine=\"470\"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n        $r.InsertXML($xml)\n        break\n    }\n}\n"}
